# Update the "Team Month" sheet data with refreshed figures, and adjust the
# selected cell to reflect where the user left off (D2) after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Team Month")

# --- Shared-string labels for rows 5 & 6 swap (B5/B6 text exchange places) ---
$ws.Range("B5").Value = "Xinghao_s1l"
$ws.Range("B6").Value = "Anisa_s1"

# --- Refreshed "Weekly Pending Total(Rp)" values ---
$ws.Range("C5").Value = 6374400927
$ws.Range("C6").Value = 7196293897

# --- Refreshed "Repayment" values (column D) ---
$ws.Range("D2").Value = 495335588
$ws.Range("D3").Value = 1843140130
$ws.Range("D4").Value = 984226650
$ws.Range("D5").Value = 1397429537
$ws.Range("D6").Value = 1574723547

# Recovery rate (column E) is a live formula (D/C) already present on the
# sheet, so it recalculates automatically from the new C/D values above.

# Reflect the user's last selection before saving.
$ws.Range("D2").Select()
